# Implementacion parcial de generar recibo de pago de alumno.
# Se actualiza la plantilla de "Lista de Tareas de la 7ma Iteración":
# se registra 1 hora consumida en el Día 2 de la tarea
# "CU Generar recibo de pago." (fila 6, hoja "Casos de Uso"),
# lo que recalcula en cascada las columnas de horas restantes/consumidas
# de los días subsiguientes y los totales de la fila.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

# Día 2 - Cons. (horas consumidas) para la tarea "CU Generar recibo de pago."
$ws.Range("K6").Value = 1

# Vuelve a aplicar las combinaciones de celdas del encabezado (fila 4) para
# las columnas de días, preservando el mismo layout.
$headerMerges = @("AL4:AM4","H4:I4","K4:L4","N4:O4","Q4:R4","T4:U4","W4:X4","Z4:AA4","AC4:AD4","AF4:AG4","AI4:AJ4")
foreach ($ref in $headerMerges) {
    $ws.Range($ref).UnMerge()
    $ws.Range($ref).Merge()
}

# Actualiza la selección activa de la hoja hacia la celda recién editada.
$ws.Range("F7").Select()
